$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.812.69"
$ws.Range("E2").Value = "  +0.82%  "

$ws.Range("D3").Value = "1.894.18"
$ws.Range("E3").Value = "  +0.48%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9998"
$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "249.73"
$ws.Range("E5").Value = "  +1.02%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9998"
$ws.Range("E6").Value = "  -0.05%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4767"
$ws.Range("E7").Value = "  +0.01%  "

$ws.Range("E8").Value = "  +0.70%  "

$ws.Range("E10").Value = "  +0.20%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07774"
$ws.Range("E11").Value = "  +0.68%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "97.62"
$ws.Range("E12").Value = "  -0.45%  "

$ws.Range("D13").Value = "1.891.23"
$ws.Range("E13").Value = "  +0.42%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.7400"
$ws.Range("E14").Value = "  -0.33%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.254"
$ws.Range("E15").Value = "  +1.90%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "284.50"
$ws.Range("E16").Value = "  +3.43%  "

$ws.Range("D17").Value = "30.871.22"
$ws.Range("E17").Value = "  +1.01%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.25"
$ws.Range("E18").Value = "  -1.94%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007593"
$ws.Range("E19").Value = "  +0.15%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9998"
$ws.Range("E20").Value = "  -0.05%  "

$ws.Range("D21").Value = "2.139.07"
$ws.Range("E21").Value = "  +0.45%  "

$ws.Range("E22").Value = "  +1.62%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.000"
$ws.Range("E23").Value = "  -0.06%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.274"
$ws.Range("E24").Value = "  +1.16%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.279"
$ws.Range("E25").Value = "  -0.75%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "164.52"
$ws.Range("E26").Value = "  +0.55%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.02"
$ws.Range("E27").Value = "  +0.28%  "

$ws.Range("E28").Value = "  -0.95%  "

$ws.Range("E29").Value = "  -1.80%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09787"
$ws.Range("E30").Value = "  -2.84%  "

$ws.Range("E31").Value = "  -0.82%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.331"
$ws.Range("E32").Value = "  -0.01%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.211"
$ws.Range("E33").Value = "  +2.01%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04907"
$ws.Range("E34").Value = "  +1.87%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.134"
$ws.Range("E35").Value = "  +0.06%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7039"
$ws.Range("E36").Value = "  +0.21%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.720"
$ws.Range("E37").Value = "  +0.17%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01917"
$ws.Range("E38").Value = "  +2.22%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.813"
$ws.Range("E39").Value = "  +2.25%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.382"
$ws.Range("E40").Value = "  +0.70%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "76.12"
$ws.Range("E41").Value = "  +6.14%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.037"
$ws.Range("E42").Value = "  +2.12%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.4289"
$ws.Range("E43").Value = "  +1.30%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9998"
$ws.Range("E44").Value = "  -0.03%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.8382"
$ws.Range("E45").Value = "  -0.40%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "102.02"
$ws.Range("E46").Value = "  -0.79%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.448"
$ws.Range("E47").Value = "  +0.94%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.107"
$ws.Range("E48").Value = "  -0.17%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "35.88"
$ws.Range("E49").Value = "  +0.42%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "929.99"
$ws.Range("E50").Value = "  +1.18%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05767"
